# Github Auto Build at 2023-12-09 14:02
# Append the latest cost-log rows (144-146) to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A144").Value = "2023-12-09 14:01:51"
$ws.Range("B144").Value = 0.0008

$ws.Range("A145").Value = "2023-12-09 14:02:00"
$ws.Range("B145").Value = 0.0006000000000000001

$ws.Range("A146").Value = "2023-12-09 14:02:03"
$ws.Range("B146").Value = 0.0004
